# Auto-applies crypto price/volume update for Sun May 21 15:19:34 UTC 2023 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.185.37"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.89"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.44"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07412"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8813"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.34"
$ws.Range("E12").Value = "  +3.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07711"
$ws.Range("E13").Value = "  +5.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.391"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.35"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.542"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008733"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.584.66"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.250"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.089.05"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.881"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.98"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.52"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.121"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.190"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.65"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08934"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7458"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.166"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.964"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.518"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.011"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.552"
$ws.Range("E37").Value = "  +6.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.094"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05300"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01940"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.306"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.947"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5279"
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.381"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4909"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.39"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.011"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.62"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.657"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06280"
$ws.Range("E51").Value = "  -0.42%  "
